# Reorder the player roster rows (A2:C18) so that each player's row
# (Name, Position, Team) is moved as a unit to its new sorted position.
# The underlying set of rows is unchanged; only their order changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Damian Lillard",      "PG",       "Milwaukee Bucks"),
    @("Derrick White",       "PG,SG",    "Boston Celtics"),
    @("Isaiah Collier",      "PG",       "Utah Jazz"),
    @("Collin Sexton",       "PG,SG",    "Utah Jazz"),
    @("Cameron Johnson",     "SF,PF",    "Brooklyn Nets"),
    @("Onyeka Okongwu",      "PF,C",     "Atlanta Hawks"),
    @("Malik Monk",          "PG,SG,SF", "Sacramento Kings"),
    @("Julius Randle",       "PF,C",     "Minnesota Timberwolves"),
    @("Naz Reid",            "PF,C",     "Minnesota Timberwolves"),
    @("Cade Cunningham",     "PG,SG",    "Detroit Pistons"),
    @("LaMelo Ball",         "PG,SG",    "Charlotte Hornets"),
    @("Anthony Davis",       "PF,C",     "Los Angeles Lakers"),
    @("Duncan Robinson",     "SG,SF",    "Miami Heat"),
    @("Deandre Ayton",       "C",        "Portland Trail Blazers"),
    @("Devin Vassell",       "SG,SF",    "San Antonio Spurs"),
    @("Isaiah Hartenstein",  "C",        "Oklahoma City Thunder"),
    @("Coby White",          "PG,SG",    "Chicago Bulls")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $row++
}
